# ATF_fsh_length_comp.xlsx — add a "sample_size" sheet (females!B + males!B,
# with matching year column pulled from females!A), and update the sheet
# selections to match the new layout.

$wb = $excel.ActiveWorkbook

$females = $wb.Worksheets.Item("females")
$males   = $wb.Worksheets.Item("males")

# --- females sheet: update selection ---
$females.Activate()
$females.Range("C1:AB39").Select()

# --- males sheet: update selection ---
$males.Activate()
$males.Range("C1:AB39").Select()

# --- add the new "sample_size" sheet after "males" (last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sampleSize = $wb.Worksheets.Add($null, $lastSheet)
$sampleSize.Name = "sample_size"

for ($r = 1; $r -le 39; $r++) {
    $sampleSize.Range("A$r").Formula = "=females!A$r"
    $sampleSize.Range("B$r").Formula = "=females!B$r+males!B$r"
}

# sample_size becomes the active sheet/tab, with F40 selected
$sampleSize.Range("F40").Select()
